# Atualização de bases das ligas, do dia: 15-04-2024 às 22:35
#
# This script:
#   1) Swaps the data (columns B..AC) of five pairs of existing rows
#      (the index/number in column A stays put - only the match data moves).
#   2) Appends two brand-new match rows (142 and 143) at the end of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($r1, $r2) {
    # Swap columns B (2) through AC (29) between two rows, leaving column A alone.
    for ($c = 2; $c -le 29; $c++) {
        $cell1 = $ws.Cells.Item($r1, $c)
        $cell2 = $ws.Cells.Item($r2, $c)
        $v1 = $cell1.Value2
        $v2 = $cell2.Value2
        $cell1.Value = $v2
        $cell2.Value = $v1
    }
}

Swap-RowData 26 27
Swap-RowData 89 90
Swap-RowData 100 102
Swap-RowData 101 104
Swap-RowData 117 118

function Set-IndexCell($row, $value) {
    # Column A style: bold, centered, bordered (matches the other data rows).
    # Copy the existing formatting from row 141 (xlPasteFormats = -4122) so we
    # reuse the existing style entry instead of growing styles.xml.
    $c = $ws.Cells.Item($row, 1)
    $c.Value = $value
    $ws.Range("A141").Copy() | Out-Null
    $ws.Cells.Item($row, 1).PasteSpecial(-4122) | Out-Null
}

function Set-DateCell($row, $value) {
    $c = $ws.Cells.Item($row, 5)
    $c.Value = $value
    $ws.Range("E141").Copy() | Out-Null
    $ws.Cells.Item($row, 5).PasteSpecial(-4122) | Out-Null
}

# ---- New row 142 ----
Set-IndexCell 142 140
$ws.Cells.Item(142, 2).Value = 7862925
$ws.Cells.Item(142, 3).Value = "Lithuania A Lyga"
$ws.Cells.Item(142, 4).Value = "Lithuania A Lyga"
Set-DateCell 142 45396.29166666666
$ws.Cells.Item(142, 6).Value = "FK Dziugas Telsiai"
$ws.Cells.Item(142, 7).Value = "FK Siauliai"
$ws.Cells.Item(142, 8).Value = 2
$ws.Cells.Item(142, 9).Value = 1
$ws.Cells.Item(142, 10).Value = "H"
$ws.Cells.Item(142, 11).Value = 3.75
$ws.Cells.Item(142, 12).Value = 3.1
$ws.Cells.Item(142, 13).Value = 1.909
$ws.Cells.Item(142, 14).Value = 4
$ws.Cells.Item(142, 15).Value = 3.2
$ws.Cells.Item(142, 16).Value = 1.85
$ws.Cells.Item(142, 17).Value = 0.5
$ws.Cells.Item(142, 18).Value = 1.875
$ws.Cells.Item(142, 19).Value = 1.925
$ws.Cells.Item(142, 20).Value = 2
$ws.Cells.Item(142, 21).Value = 1.8
$ws.Cells.Item(142, 22).Value = 2
$ws.Cells.Item(142, 23).Value = 3
$ws.Cells.Item(142, 24).Value = -1
$ws.Cells.Item(142, 25).Value = -1
$ws.Cells.Item(142, 26).Value = 0.875
$ws.Cells.Item(142, 27).Value = -1
$ws.Cells.Item(142, 28).Value = 0.8
$ws.Cells.Item(142, 29).Value = -1

# ---- New row 143 ----
Set-IndexCell 143 141
$ws.Cells.Item(143, 2).Value = 7862046
$ws.Cells.Item(143, 3).Value = "Lithuania A Lyga"
$ws.Cells.Item(143, 4).Value = "Lithuania A Lyga"
Set-DateCell 143 45396.375
$ws.Cells.Item(143, 6).Value = "Panevezys"
$ws.Cells.Item(143, 7).Value = "FK Zalgiris Vilnius"
$ws.Cells.Item(143, 8).Value = 1
$ws.Cells.Item(143, 9).Value = 2
$ws.Cells.Item(143, 10).Value = "A"
$ws.Cells.Item(143, 11).Value = 3.75
$ws.Cells.Item(143, 12).Value = 3.25
$ws.Cells.Item(143, 13).Value = 1.833
$ws.Cells.Item(143, 14).Value = 3.5
$ws.Cells.Item(143, 15).Value = 3
$ws.Cells.Item(143, 16).Value = 2
$ws.Cells.Item(143, 17).Value = 0.25
$ws.Cells.Item(143, 18).Value = 2.025
$ws.Cells.Item(143, 19).Value = 1.775
$ws.Cells.Item(143, 20).Value = 2
$ws.Cells.Item(143, 21).Value = 1.825
$ws.Cells.Item(143, 22).Value = 1.975
$ws.Cells.Item(143, 23).Value = -1
$ws.Cells.Item(143, 24).Value = -1
$ws.Cells.Item(143, 25).Value = 1
$ws.Cells.Item(143, 26).Value = -1
$ws.Cells.Item(143, 27).Value = 0.7749999999999999
$ws.Cells.Item(143, 28).Value = 0.825
$ws.Cells.Item(143, 29).Value = -1
